$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D whose new values look numeric need an explicit text
# number format first, otherwise Excel auto-converts them to numbers
# (e.g. "1.00" -> 1, "0.370" -> 0.37, "0.0000178" -> 1.78E-05).
$ws.Range('D2').Value = '63.850.84'
$ws.Range('E2').Value = '  +1.46%  '
$ws.Range('D3').Value = '2.518.40'
$ws.Range('E3').Value = '  +2.22%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '578.05'
$ws.Range('E5').Value = '  +0.41%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '151.07'
$ws.Range('E6').Value = '  +3.07%  '
$ws.Range('E7').Value = '  +0.04%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.536'
$ws.Range('E8').Value = '  -0.12%  '
$ws.Range('E9').Value = '  +0.01%  '
$ws.Range('E10').Value = '  -0.80%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.25'
$ws.Range('E11').Value = '  -0.49%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.351'
$ws.Range('E12').Value = '  -1.54%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '29.35'
$ws.Range('E13').Value = '  +1.04%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.0000178'
$ws.Range('E14').Value = '  -0.27%  '
$ws.Range('D15').Value = '2.978.14'
$ws.Range('E15').Value = '  +2.30%  '
$ws.Range('D16').Value = '63.820.27'
$ws.Range('E16').Value = '  +1.45%  '
$ws.Range('D17').Value = '2.527.54'
$ws.Range('E17').Value = '  +2.34%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '7.76'
$ws.Range('E18').Value = '  -2.60%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '10.88'
$ws.Range('E19').Value = '  -1.45%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '4.23'
$ws.Range('E20').Value = '  +2.21%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '325.53'
$ws.Range('E21').Value = '  -0.62%  '
$ws.Range('E22').Value = '  -0.44%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.00'
$ws.Range('E23').Value = '  +0.13%  '
$ws.Range('E24').Value = '  -3.27%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '65.39'
$ws.Range('E25').Value = '  -0.73%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '655.46'
$ws.Range('E26').Value = '  +0.52%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.0000102'
$ws.Range('E27').Value = '  +3.46%  '
$ws.Range('D28').Value = '2.655.69'
$ws.Range('E28').Value = '  +2.78%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.997'
$ws.Range('E29').Value = '  -0.33%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.46'
$ws.Range('E30').Value = '  +0.33%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.99'
$ws.Range('E31').Value = '  -0.34%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.84'
$ws.Range('E32').Value = '  -0.98%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.134'
$ws.Range('E33').Value = '  +0.65%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.998'
$ws.Range('E34').Value = '  -0.01%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.52'
$ws.Range('E35').Value = '  -1.57%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '4.75'
$ws.Range('E36').Value = '  -0.24%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.49'
$ws.Range('E37').Value = '  +1.08%  '
$ws.Range('B38').Value = 'PolygonEcosystemToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.370'
$ws.Range('E38').Value = '  -0.07%  '
$ws.Range('B39').Value = 'Monero'
$ws.Range('C39').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '151.92'
$ws.Range('E39').Value = '  -0.29%  '
$ws.Range('B40').Value = 'dogwifhat'
$ws.Range('C40').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.80'
$ws.Range('E40').Value = '  +0.88%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '18.72'
$ws.Range('E41').Value = '  -0.12%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.76'
$ws.Range('E42').Value = '  +1.43%  '
$ws.Range('E43').Value = '  +0.00%  '
$ws.Range('B44').Value = 'Aave'
$ws.Range('C44').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '157.72'
$ws.Range('E44').Value = '  +2.51%  '
$ws.Range('B45').Value = 'BabyDogeCoin'
$ws.Range('C45').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D45').Value = '0.0₆0301'
$ws.Range('E45').Value = '  -5.30%  '
$ws.Range('E46').Value = '  +1.22%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.62'
$ws.Range('E47').Value = '  +0.71%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '20.71'
$ws.Range('E48').Value = '  +1.07%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.613'
$ws.Range('E49').Value = '  +0.78%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0516'
$ws.Range('E50').Value = '  +0.53%  '
$ws.Range('B51').Value = 'Stellar'
$ws.Range('C51').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0911'
$ws.Range('E51').Value = '  -0.54%  '
